# Update gh-pages to output generated at 456a3b4
# Applies the scraped-data refresh to the 广州-漫展信息 workbook:
#  - "想去人数" (interest count) bumps in column F across the 展览 / 演出 /
#    本地生活 sheets and their mirrored rows in 全部类型
#  - one changed cover-image URL
#  - two events whose 最低票价 (column G) flipped from a numeric price to
#    the literal status text "已售罄" (sold out)

$wb = $excel.ActiveWorkbook

$ws_exhibit = $wb.Worksheets.Item("展览")
$ws_show    = $wb.Worksheets.Item("演出")
$ws_local   = $wb.Worksheets.Item("本地生活")
$ws_all     = $wb.Worksheets.Item("全部类型")

# --- 展览 (exhibitions) sheet, plus its mirrored rows in 全部类型 ---

$ws_exhibit.Range("F2").Value = 8206
$ws_all.Range("F3").Value = 8206

$ws_exhibit.Range("F5").Value = 33009
$ws_all.Range("F7").Value = 33010

$ws_exhibit.Range("F14").Value = 625
$ws_all.Range("F20").Value = 625

$ws_exhibit.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202403/asJLaLl91711679266704.jpeg"
$ws_all.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202403/asJLaLl91711679266704.jpeg"

$ws_exhibit.Range("F15").Value = 424
$ws_all.Range("F21").Value = 424

$ws_exhibit.Range("F17").Value = 569
$ws_all.Range("F28").Value = 569

$ws_exhibit.Range("F19").Value = 426
$ws_all.Range("F30").Value = 426

$ws_exhibit.Range("F23").Value = 729
$ws_all.Range("F34").Value = 729

$ws_exhibit.Range("F24").Value = 2378
$ws_all.Range("F35").Value = 2378

$ws_exhibit.Range("F25").Value = 855
$ws_all.Range("F36").Value = 855

$ws_exhibit.Range("F31").Value = 8
$ws_all.Range("F43").Value = 8

$ws_exhibit.Range("F32").Value = 1087
$ws_all.Range("F44").Value = 1087

# --- 演出 (shows) sheet, plus its mirrored rows in 全部类型 ---

$ws_show.Range("G3").Value = "已售罄"
$ws_all.Range("G12").Value = "已售罄"

$ws_show.Range("F8").Value = 32
$ws_all.Range("F24").Value = 32

# --- 本地生活 (local life) sheet, plus its mirrored row in 全部类型 ---

$ws_local.Range("F2").Value = 539
$ws_all.Range("F2").Value = 539
